$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Rename column J ("מס' טלפון" -> "טלפון אזרחי") and add a new trailing column ("סטטוס").
$ws.Range("J1").Value = "טלפון אזרחי"
$newCol = $tbl.ListColumns.Add()
$ws.Range("L1").Value = "סטטוס"

# Header row formatting: bold, centered/top aligned, thin-bordered.
$headerRange = $tbl.HeaderRowRange
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$ws.Range("E16").Select()
